$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper approach: for Price (column D) values that look numeric (e.g. "206.02"),
# Excel COM auto-converts the assigned string into a real number, which would
# strip the original text formatting (e.g. trailing zeros, multi-dot grouping).
# These prices must stay plain text, so force a temporary text NumberFormat,
# assign the value, then restore the default "Normal" style so no extra
# explicit cell style is left behind (matching the source cells' lack of 's' attr).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.854.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.563.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.74"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.34%  "
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("E10").Value = "  -1.30%  "
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.786.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.573.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.861.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("E19").Value = "  +2.00%  "
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.32%  "
$ws.Range("E24").Value = "  +1.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  -1.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0466"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.62%  "
$ws.Range("E31").Value = "  -3.61%  "
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.403.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.55%  "
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("E37").Value = "  -2.94%  "
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.529"
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.995"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.699.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("E49").Value = "  +2.71%  "
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("E51").Value = "  +0.61%  "
